# edit.ps1
# Applies the "jk added some work" commit:
#   1. Marks the three inline-picture runs as NoProofing (adds <w:rPr><w:noProof/></w:rPr>
#      to the <w:r> that hosts each <w:drawing>).
#   2. Changes the hard-coded loop bound `6` in `for(int i=0; i<6; i++){` to the
#      identifier `num_cities`, matching the surrounding Consolas/E06C75 code-run styling
#      and splitting the edit into two runs (" " then "num_cities") like the authored
#      source.

$d = $word.ActiveDocument

# --- 1. Tag every inline picture's run with <w:noProof/> -------------------
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $shape.Range.NoProofing = $true
}

# --- 2. Replace the loop bound "6" with "num_cities" ------------------------
# Locate the lone "6" inside "i<6; i++" (the other standalone "6" earlier in the
# doc belongs to the "ENSF 614" course number, so search forward and inspect).
$rng = $d.Content.Duplicate
$rng.Find.Text = "6"
$rng.Find.Forward = $true
$rng.Find.Wrap = 0
$rng.Find.MatchCase = $true
$rng.Find.MatchWholeWord = $false

$target = $null
while ($rng.Find.Execute()) {
    if ($rng.Font.Color -eq 6724305) {
        $target = $d.Range($rng.Start, $rng.End)
        break
    }
    $rng.Collapse(0)
}

if ($target -ne $null) {
    # Turn the "6" run into a single space, recoloured to match the other
    # identifiers on this code line (E06C75), then append "num_cities" right
    # after it as its own run (forcing the run split via a transient Bold
    # toggle, then clearing it again so the final formatting matches).
    $target.Text = " "
    $target.Font.Color = 7695584
    $target.Collapse(0)
    $target.InsertAfter("num_cities")
    $target.Bold = 1
    $target.Bold = 0
}

Write-Output "done"
